$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 26 and 27 with new data (weekly roll-forward)
$ws.Range("D26").Value = 45212
$ws.Range("D27").Value = 45194
$ws.Range("M27").Value = 80

# Append new row 28, which is a copy of the former row 27 data (before its update)
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 45176
$ws.Range("D28").NumberFormat = $ws.Range("D27").NumberFormat
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = "Otros"
$ws.Range("I28").Value = 100107002
$ws.Range("J28").Value = "Chirimoya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 30
$ws.Range("N28").Value = 22000
$ws.Range("O28").Value = 22000
$ws.Range("P28").Value = 22000
$ws.Range("Q28").Value = "$/bandeja 10 kilos"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 2200
$ws.Range("T28").Value = 10
